$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 135; this shifts the existing rows 135-146 down to 136-147,
# carrying the date-format style on column D along with them.
$ws.Rows.Item(135).Insert()

# Populate the newly inserted row 135 with the new weekly data point.
$ws.Range("A135").Value = 10
$ws.Range("B135").Value = "Vega Modelo de Temuco"
$ws.Range("C135").Value = "La Araucanía"
$ws.Range("D135").Value = 44461
$ws.Range("E135").Value = 9
$ws.Range("F135").Value = 100112039
$ws.Range("G135").Value = "Ciboulette"
$ws.Range("H135").Value = "Sin especificar"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 40
$ws.Range("K135").Value = 3000
$ws.Range("L135").Value = 6000
$ws.Range("M135").Value = 3750
$ws.Range("N135").Value = "$/docena de atados"
$ws.Range("O135").Value = "Provincia de Cautín"
$ws.Range("P135").Value = 1250
$ws.Range("Q135").Value = 3
$ws.Range("R135").Value = "Hortaliza"
